$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "64.380.54"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -2.47%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.178.85"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -4.01%  "

# Row 4
$ws.Range("E4").Value = "  +0.01%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "570.86"
$ws.Range("D5").Style = "Normal"

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "168.75"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -7.44%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.608"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -6.00%  "

# Row 8
$ws.Range("E8").Value = "  -0.14%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "3.187.70"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.67%  "

# Row 10
$ws.Range("E10").Value = "  -3.52%  "

# Row 11
$ws.Range("E11").Value = "  -0.34%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.387"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -3.06%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.730.99"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -4.08%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.130"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.59%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "64.469.32"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.48%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "25.30"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -3.04%  "

# Row 17
$ws.Range("E17").Value = "  -3.57%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.190.87"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -4.00%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "418.21"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.31%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.96"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.14%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.36"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.12%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.12"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.30%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.00"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.01%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "70.36"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.87%  "

# Row 25
$ws.Range("E25").Value = "  -0.02%  "

# Row 26
$ws.Range("E26").Value = "  +2.96%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.489"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -4.27%  "

# Row 28
$ws.Range("E28").Value = "  -6.14%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.76"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.47%  "

# Row 30
$ws.Range("E30").Value = "  +0.01%  "

# Row 31
$ws.Range("E31").Value = "  -3.45%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "21.75"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.53%  "

# Row 33
$ws.Range("E33").Value = "  -0.08%  "

# Row 34
$ws.Range("E34").Value = "  -1.88%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.35"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.83%  "

# Row 36
$ws.Range("E36").Value = "  -3.54%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "156.93"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.36%  "

# Row 38
$ws.Range("E38").Value = "  -4.62%  "

# Row 39
$ws.Range("B39").Value = "Stacks"
$ws.Range("C39").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.70"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -4.69%  "

# Row 40
$ws.Range("B40").Value = "Maker"
$ws.Range("C40").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.696.03"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -6.17%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "24.25"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -7.60%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.21"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.88%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "39.24"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.70%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.717"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -5.53%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0623"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -5.43%  "

# Row 46
$ws.Range("E46").Value = "  -5.33%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0264"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.35%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "292.55"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -6.45%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "21.48"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -6.68%  "

# Row 50
$ws.Range("E50").Value = "  -11.35%  "

# Row 51
$ws.Range("E51").Value = "  -0.16%  "
